# Add a PROFIT column (K) to the SalesData sheet, computing a per-product-line
# margin off the SALES (H) column, and fix one mis-categorised row
# (row 77: HDD25-1TB -> HDD25-1TB-SG1T, which is really EXTERNALSTORAGE, not
# INTERNALSTORAGE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 135

# --- 1. Fix row 77's product line / code -----------------------------------
# It was miscoded as INTERNALSTORAGE / HDD25-1TB; it should be an external
# drive, coded EXTERNALSTORAGE / HDD25-1TB-SG1T.
$ws.Range("E77").Value = "EXTERNALSTORAGE"
$ws.Range("F77").Value = "HDD25-1TB-SG1T"

# --- 2. Add the PROFIT header in K1 -----------------------------------------
$ws.Range("K1").Value = "PROFIT"

# --- 3. Fill K2:K135 with the profit formula --------------------------------
# PROFIT = SALES(H) scaled by a margin percentage that depends on the
# PRODUCTLINE (E), with EXTERNALSTORAGE sales booked at a flat 100000.
for ($r = 2; $r -le $lastRow; $r++) {
    $formula = '=IF(E' + $r + '="INTERNALSTORAGE",H' + $r + '*52%,' +
        'IF(E' + $r + '="SERVICE",H' + $r + '*100%,' +
        'IF(E' + $r + '="DISPLAY",H' + $r + '*50%,' +
        'IF(E' + $r + '="RAM",H' + $r + '*65%,' +
        'IF(E' + $r + '="SOFTWARE",H' + $r + '*100%,' +
        'IF(E' + $r + '="BATTERY",H' + $r + '*70%,' +
        'IF(E' + $r + '="KEYBOARD",H' + $r + '*60%,' +
        'IF(E' + $r + '="ADAPTOR",H' + $r + '*45%,' +
        'IF(E' + $r + '="MAINBOARD",H' + $r + '*15%,' +
        'IF(E' + $r + '="ACCESSORIES",H' + $r + '*25%,' +
        'IF(E' + $r + '="VGA",H' + $r + '*15%,' +
        'IF(E' + $r + '="POWERSUPPLY",H' + $r + '*40%,' +
        'IF(E' + $r + '="PROCESSOR",H' + $r + '*20%,' +
        'IF(E' + $r + '="SECOND",H' + $r + '*150%,' +
        'IF(E' + $r + '="CASING",H' + $r + '*15%,' +
        'IF(E' + $r + '="MONITOR",H' + $r + '*5%,' +
        'IF(E' + $r + '="EXTERNALSTORAGE",100000,H' + $r + ')))))))))))))))))'
    $ws.Range("K$r").Formula = $formula
}

# --- 4. Extend the table range to include the new column -------------------
$ws.Range("A1:K$lastRow").AutoFilter()

$filterName = $wb.Names.Item("_xlnm._FilterDatabase")
$filterName.RefersTo = "=SalesData!`$A`$1:`$K`$" + $lastRow

# --- 5. Move the selection to the new column so it is visible --------------
$ws.Range("K1").Select()
